$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Merge the three header rows (A2:F2, A3:F3, A4:F4) into a single
# --- merged block A2:F4 holding the combined template placeholder text.

# Start from a clean (unmerged) state so the three original merged ranges
# can be re-merged as one.
$ws.Range("A2:F2").UnMerge()
$ws.Range("A3:F3").UnMerge()
$ws.Range("A4:F4").UnMerge()

# Drop the borders that used to decorate these header cells - the new
# merged block carries no border - and force the "Normal 100" cell style
# (same named style already used throughout the sheet) so Excel reuses the
# existing style record instead of minting one off the default "Normal"
# style.
$hdrRange = $ws.Range("A2:F4")
$hdrRange.Borders.LineStyle = -4142
$hdrRange.Style = "Normal 100"

# Combine the three separate placeholders into one multi-line value living
# in A2; clear out what used to live in A3 / A4.
$ws.Range("A2").Value2 = "{{NombreReporte}}`n{{Detalle}}`n{{Compania}}"
$ws.Range("A3").Value2 = ""
$ws.Range("A4").Value2 = ""

# Re-merge as a single block and apply the centered / wrapped formatting.
$hdrRange.Merge()
$hdrRange.HorizontalAlignment = -4108
$hdrRange.VerticalAlignment = -4108
$hdrRange.WrapText = $true

# Keep the original row heights for rows 2-4 (merging/wrapping would
# otherwise trigger auto-fit and change them).
$ws.Rows(2).RowHeight = 17.25
$ws.Rows(3).RowHeight = 15
$ws.Rows(4).RowHeight = 15.75

# Row 9 (the item placeholder row) gets a taller, explicit row height.
$ws.Rows(9).RowHeight = 30

# Move the active selection, matching the author's last cursor position.
$ws.Range("D14").Select()
